$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New release column: 11.4.0 (column I) ---

# Header (I1), bold like the rest of row 1 (copy header style from H1)
$ws.Cells.Item(1, 8).Copy()
$ws.Cells.Item(1, 9).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Cells.Item(1, 9).Value = "11.4.0"

# New package version for Tardigrade.Framework.EntityFramework at 11.4.0
$ws.Cells.Item(7, 9).Value = "10.0.0"

# The remaining rows in column I (3,4,6,8,9,10,11) stay blank for this
# release, but keep a touched/formatted cell in place (matching column H's
# pattern of blank-but-styled cells) so the column is fully populated.
$blankRows = @(3, 4, 6, 8, 9, 10, 11)
foreach ($r in $blankRows) {
    $ws.Cells.Item($r, 9).Borders.LineStyle = -4142
}

# Move the active selection to I7, where the new value was entered
$ws.Range("I7").Select() | Out-Null
